$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two team names (B6/B7) to match new shared-string ordering
$ws.Range("B6").Value = "Hansyah_S2l"
$ws.Range("B7").Value = "EDN_S2l"

# Update Weekly Pending Total(Rp) and Repayment figures (column E recomputes via existing formula)
$ws.Range("C2").Value = 7442439054
$ws.Range("D2").Value = 887111589

$ws.Range("C3").Value = 3081566235
$ws.Range("D3").Value = 350470555

$ws.Range("C4").Value = 3995841125
$ws.Range("D4").Value = 437178785

$ws.Range("C5").Value = 2042479922
$ws.Range("D5").Value = 218250286

$ws.Range("C6").Value = 7531470001
$ws.Range("D6").Value = 795868021

$ws.Range("C7").Value = 2047822376
$ws.Range("D7").Value = 211943320

$ws.Range("C8").Value = 7374789257
$ws.Range("D8").Value = 748361716

# Move/collapse the active selection to D2
$ws.Range("D2").Select()
